# Generate Report for Handback
# ------------------------------------------------------------------
# This script mirrors a "handback" report-generation pass over the
# localization-status workbook:
#   1. The overall status text "Ready for handoff" is updated, in every
#      cell that shows it, to "Handed back: in sync with en-US".
#   2. The per-language sheets (zh-cn, de-de) get their "Latest Target
#      File" / "Latest Handback File" columns (F/G) populated with the
#      handed-back source/file names, each one a hyperlink, for both
#      data rows.
#   3. The de-de sheet's "Latest Handback DateTime" column (H) is
#      stamped with the actual handback timestamp for both rows.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- 1. Update every cell currently showing "Ready for handoff" -----
# (Overview status columns, plus the Status column on each language
# sheet - they all share the same underlying string, so they must all
# be updated together.)
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# --- 2. Populate Latest Target File (F) / Latest Handback File (G) --

# zh-cn, row 2 (05f62907-...)
$zhcn.Hyperlinks.Add(
    $zhcn.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/ccde773797c59012d3b4bd6458042d5fa86cb703/e2e/05f62907-ca4a-4e1a-92c1-1e5e9375ded2.md",
    "",
    "",
    "05f62907-ca4a-4e1a-92c1-1e5e9375ded2.md"
)
$zhcn.Hyperlinks.Add(
    $zhcn.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/ccde773797c59012d3b4bd6458042d5fa86cb703/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/05f62907-ca4a-4e1a-92c1-1e5e9375ded2.4aa060e253e9e8673f2a45f0f7cde34898b76145.zh-cn.xlf",
    "",
    "",
    "05f62907-ca4a-4e1a-92c1-1e5e9375ded2.4aa060e253e9e8673f2a45f0f7cde34898b76145.zh-cn.xlf"
)

# zh-cn, row 3 (69f72130-...)
$zhcn.Hyperlinks.Add(
    $zhcn.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/ccde773797c59012d3b4bd6458042d5fa86cb703/e2e/69f72130-834e-4bac-a39b-243b12f1ca4b.md",
    "",
    "",
    "69f72130-834e-4bac-a39b-243b12f1ca4b.md"
)
$zhcn.Hyperlinks.Add(
    $zhcn.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/ccde773797c59012d3b4bd6458042d5fa86cb703/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/high/69f72130-834e-4bac-a39b-243b12f1ca4b.fe7ddf077176cbdd789dec0a497278754d466e15.zh-cn.xlf",
    "",
    "",
    "69f72130-834e-4bac-a39b-243b12f1ca4b.fe7ddf077176cbdd789dec0a497278754d466e15.zh-cn.xlf"
)

# de-de, row 2 (05f62907-...)
$dede.Hyperlinks.Add(
    $dede.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/ccde773797c59012d3b4bd6458042d5fa86cb703/e2e/05f62907-ca4a-4e1a-92c1-1e5e9375ded2.md",
    "",
    "",
    "05f62907-ca4a-4e1a-92c1-1e5e9375ded2.md"
)
$dede.Hyperlinks.Add(
    $dede.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/090dd8ef1d25519530fea748e721e72023833456/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/high/05f62907-ca4a-4e1a-92c1-1e5e9375ded2.4aa060e253e9e8673f2a45f0f7cde34898b76145.de-de.xlf",
    "",
    "",
    "05f62907-ca4a-4e1a-92c1-1e5e9375ded2.4aa060e253e9e8673f2a45f0f7cde34898b76145.de-de.xlf"
)

# de-de, row 3 (69f72130-...)
$dede.Hyperlinks.Add(
    $dede.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/ccde773797c59012d3b4bd6458042d5fa86cb703/e2e/69f72130-834e-4bac-a39b-243b12f1ca4b.md",
    "",
    "",
    "69f72130-834e-4bac-a39b-243b12f1ca4b.md"
)
$dede.Hyperlinks.Add(
    $dede.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/090dd8ef1d25519530fea748e721e72023833456/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/high/69f72130-834e-4bac-a39b-243b12f1ca4b.fe7ddf077176cbdd789dec0a497278754d466e15.de-de.xlf",
    "",
    "",
    "69f72130-834e-4bac-a39b-243b12f1ca4b.fe7ddf077176cbdd789dec0a497278754d466e15.de-de.xlf"
)

# --- 3. Stamp the Latest Handback DateTime (H) for de-de ------------
$dede.Range("H2").Value = "2016-03-22 00:14:22"
$dede.Range("H3").Value = "2016-03-22 00:14:22"
